$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows before row 334, pushing existing rows 334-447
# down to 336-449 (matches the weekly-update pattern: a new week's data is
# inserted near the top of the "Cilantro" price history and the rest of the
# series shifts down).
$ws.Rows("334:335").Insert()

# --- New row 334: Primera / $/caja 36 atados ---
$ws.Cells.Item(334, 1).Value = 9
$ws.Cells.Item(334, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(334, 3).Value = "Metropolitana"
$ws.Cells.Item(334, 4).Value = 44559
$ws.Cells.Item(334, 5).Value = 13
$ws.Cells.Item(334, 6).Value = 100112040
$ws.Cells.Item(334, 7).Value = "Cilantro"
$ws.Cells.Item(334, 8).Value = "Sin especificar"
$ws.Cells.Item(334, 9).Value = "Primera"
$ws.Cells.Item(334, 10).Value = 34
$ws.Cells.Item(334, 11).Value = 8000
$ws.Cells.Item(334, 12).Value = 8000
$ws.Cells.Item(334, 13).Value = 8000
$ws.Cells.Item(334, 14).Value = "`$/caja 36 atados"
$ws.Cells.Item(334, 15).Value = "Región Metropolitana"
$ws.Cells.Item(334, 16).Value = 222
$ws.Cells.Item(334, 17).Value = 36
$ws.Cells.Item(334, 18).Value = "Hortaliza"

# --- New row 335: Primera / $/docena de atados ---
$ws.Cells.Item(335, 1).Value = 9
$ws.Cells.Item(335, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(335, 3).Value = "Metropolitana"
$ws.Cells.Item(335, 4).Value = 44559
$ws.Cells.Item(335, 5).Value = 13
$ws.Cells.Item(335, 6).Value = 100112040
$ws.Cells.Item(335, 7).Value = "Cilantro"
$ws.Cells.Item(335, 8).Value = "Sin especificar"
$ws.Cells.Item(335, 9).Value = "Primera"
$ws.Cells.Item(335, 10).Value = 160
$ws.Cells.Item(335, 11).Value = 14000
$ws.Cells.Item(335, 12).Value = 16000
$ws.Cells.Item(335, 13).Value = 15000
$ws.Cells.Item(335, 14).Value = "`$/docena de atados"
$ws.Cells.Item(335, 15).Value = "Región Metropolitana"
$ws.Cells.Item(335, 16).Value = 5000
$ws.Cells.Item(335, 17).Value = 3
$ws.Cells.Item(335, 18).Value = "Hortaliza"
